$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A row 3 contains "リリー" (Lily) -> add an English label in column C
$ws.Range("C3").Value = "Rewrite     -   Lily"

# Column A row 30 contains "シィナ" (Sina) -> add an English label in column C
$ws.Range("C30").Value = "Sina"
